$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [PSCustomObject]@{ Row = 2; B = "Bitcoin"; C = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D = "68.146.75"; E = "  +0.86%  " },
    [PSCustomObject]@{ Row = 3; B = "Ethereum"; C = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D = "3.689.28"; E = "  +0.28%  " },
    [PSCustomObject]@{ Row = 4; B = "TetherUSD"; C = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D = "1.00"; E = "  +0.09%  " },
    [PSCustomObject]@{ Row = 5; B = "BNB"; C = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D = "605.14"; E = "  +4.24%  " },
    [PSCustomObject]@{ Row = 6; B = "Solana"; C = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D = "191.26"; E = "  +13.31%  " },
    [PSCustomObject]@{ Row = 7; B = "XRP"; C = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D = "0.628"; E = "  +1.30%  " },
    [PSCustomObject]@{ Row = 8; B = "USDC"; C = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D = "1.00"; E = "  +0.02%  " },
    [PSCustomObject]@{ Row = 9; B = "Cardano"; C = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D = "0.714"; E = "  +2.17%  " },
    [PSCustomObject]@{ Row = 10; B = "Avalanche"; C = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D = "58.84"; E = "  +15.72%  " },
    [PSCustomObject]@{ Row = 11; B = "Dogecoin"; C = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D = "0.157"; E = "  -1.94%  " },
    [PSCustomObject]@{ Row = 12; B = "ShibaInu"; C = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D = "0.0000281"; E = "  -1.31%  " },
    [PSCustomObject]@{ Row = 13; B = "Polkadot"; C = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D = "10.31"; E = "  -0.28%  " },
    [PSCustomObject]@{ Row = 14; B = "WrappedliquidstakedEther2.0"; C = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D = "4.270.93"; E = "  +0.15%  " },
    [PSCustomObject]@{ Row = 15; B = "WrappedEther"; C = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D = "3.685.14"; E = "  +0.09%  " },
    [PSCustomObject]@{ Row = 16; B = "TRON"; C = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D = "0.127"; E = "  +0.96%  " },
    [PSCustomObject]@{ Row = 17; B = "Chainlink"; C = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D = "19.22"; E = "  +0.23%  " },
    [PSCustomObject]@{ Row = 18; B = "Polygon"; C = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D = "1.13"; E = "  +2.36%  " },
    [PSCustomObject]@{ Row = 19; B = "Uniswap"; C = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D = "12.67"; E = "  -0.31%  " },
    [PSCustomObject]@{ Row = 20; B = "WrappedBTC"; C = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D = "67.952.79"; E = "  +0.81%  " },
    [PSCustomObject]@{ Row = 21; B = "BitcoinCash"; C = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D = "404.19"; E = "  +0.72%  " },
    [PSCustomObject]@{ Row = 22; B = "PancakeSwap"; C = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D = "4.53"; E = "  +1.14%  " },
    [PSCustomObject]@{ Row = 23; B = "Litecoin"; C = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D = "89.32"; E = "  +2.32%  " },
    [PSCustomObject]@{ Row = 24; B = "RenderToken"; C = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D = "11.71"; E = "  +7.53%  " },
    [PSCustomObject]@{ Row = 25; B = "ImmutableX"; C = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D = "3.02"; E = "  -0.24%  " },
    [PSCustomObject]@{ Row = 26; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "12.81"; E = "  +1.27%  " },
    [PSCustomObject]@{ Row = 27; B = "LEO"; C = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"; D = "6.01"; E = "  -0.21%  " },
    [PSCustomObject]@{ Row = 28; B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "3.72"; E = "  -1.47%  " },
    [PSCustomObject]@{ Row = 29; B = "Filecoin"; C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D = "9.49"; E = "  +1.38%  " },
    [PSCustomObject]@{ Row = 30; B = "EthereumClassic"; C = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D = "32.25"; E = "  -0.16%  " },
    [PSCustomObject]@{ Row = 31; B = "NEARProtocol"; C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D = "7.57"; E = "  +2.55%  " },
    [PSCustomObject]@{ Row = 32; B = "InjectiveProtocol"; C = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"; D = "46.64"; E = "  +9.01%  " },
    [PSCustomObject]@{ Row = 33; B = "Cosmos"; C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "12.53"; E = "  +1.38%  " },
    [PSCustomObject]@{ Row = 34; B = "OKB"; C = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D = "67.43"; E = "  +4.27%  " },
    [PSCustomObject]@{ Row = 35; B = "Hedera"; C = "https://coinranking.com/coin/jad286TjB+hedera-hbar"; D = "0.120"; E = "  +3.89%  " },
    [PSCustomObject]@{ Row = 36; B = "Bittensor"; C = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D = "624.76"; E = "  +4.78%  " },
    [PSCustomObject]@{ Row = 37; B = "TheGraph"; C = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"; D = "0.405"; E = "  +2.65%  " },
    [PSCustomObject]@{ Row = 38; B = "Dai"; C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.00"; E = "  +0.03%  " },
    [PSCustomObject]@{ Row = 39; B = "PEPE"; C = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"; D = "0.0₃0807"; E = "  -8.62%  " },
    [PSCustomObject]@{ Row = 40; B = "FirstDigitalUSD"; C = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; D = "0.999"; E = "  -0.15%  " },
    [PSCustomObject]@{ Row = 41; B = "Kaspa"; C = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D = "0.136"; E = "  +2.27%  " },
    [PSCustomObject]@{ Row = 42; B = "ThetaToken"; C = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"; D = "2.96"; E = "  -0.10%  " },
    [PSCustomObject]@{ Row = 43; B = "VeChain"; C = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D = "0.0437"; E = "  +1.25%  " },
    [PSCustomObject]@{ Row = 44; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "2.59"; E = "  -5.74%  " },
    [PSCustomObject]@{ Row = 45; B = "Maker"; C = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D = "2.858.64"; E = "  +3.18%  " },
    [PSCustomObject]@{ Row = 46; B = "Stellar"; C = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D = "0.138"; E = "  +4.44%  " },
    [PSCustomObject]@{ Row = 47; B = "THORChain"; C = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"; D = "9.07"; E = "  -0.85%  " },
    [PSCustomObject]@{ Row = 48; B = "Monero"; C = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D = "144.63"; E = "  +3.07%  " },
    [PSCustomObject]@{ Row = 49; B = "WEMIXToken"; C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D = "2.66"; E = "  -0.54%  " },
    [PSCustomObject]@{ Row = 50; B = "ApeXProtocol"; C = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"; D = "3.05"; E = "  -2.97%  " },
    [PSCustomObject]@{ Row = 51; B = "dogwifhat"; C = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"; D = "2.57"; E = "  -6.75%  " }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 4).Style = "Normal"
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 5).Style = "Normal"
}